# Add SnpEff publication (PMID 22728672) as a new row 86 in the
# "Table3" table on the Publications sheet, shifting existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 86 (existing rows 86.. shift down to 87..)
$ws.Rows.Item(86).Insert()

$newRow = 86

$ws.Cells.Item($newRow, 1).Value = 22728672
$ws.Cells.Item($newRow, 2).Value = "SnpEff"
$ws.Cells.Item($newRow, 3).Value = "MNV;SNV;deletion;insertion"
$ws.Cells.Item($newRow, 4).Value = "3_prime_UTR_truncation;3_prime_UTR_variant;5_prime_UTR_premature_start_codon_gain_variant;5_prime_UTR_truncation;5_prime_UTR_variant;CDS_variant;chromosome_number_variation;coding_sequence_variant;coding_sequence_variant;conservative_inframe_deletion;conservative_inframe_insertion;conserved_intergenic_variant;conserved_intron_variant;disruptive_inframe_deletion;disruptive_inframe_insertion;downstream_gene_variant;exon_loss;exon_loss_variant;exon_loss_variant;exon_variant;feature_elongation;feature_truncation;frameshift_variant;gene_variant;initiator_codon_variant;intergenic_variant;intragenic_variant;intron_variant;mature_miRNA_variant;miRNA;missense_variant;NMD_transcript_variant;non_coding_transcript_exon_variant;non_coding_transcript_variant;nonsynonymous_variant;rare_amino_acid_variant;regulatory_region_ablation;regulatory_region_amplification;regulatory_region_variant;splice_acceptor_variant;splice_donor_variant;splice_region_variant;splice_region_variant;splice_site_variant;start_gained;start_lost;start_retained;stop_gained;stop_lost;stop_retained_variant;synonymous_variant;TF_binding_site_variant;TFBS_ablation;TFBS_amplification;transcript_ablation;transcript_amplification;transcript_variant;upstream_gene_variant"
$ws.Cells.Item($newRow, 6).Value = "A program for annotating and predicting the effects of single nucleotide polymorphisms, SnpEff: SNPs in the genome of Drosophila melanogaster strain w1118; iso-2; iso-3"
$ws.Cells.Item($newRow, 7).Value = "Cingolani P, Platts A, Wang le L, Coon M, Nguyen T, Wang L, Land SJ, Lu X, Ruden DM."
$ws.Cells.Item($newRow, 8).Value = "Fly (Austin). 2012 Apr-Jun;6(2):80-92. doi: 10.4161/fly.19695."
$ws.Cells.Item($newRow, 9).Value = "Cingolani P"
$ws.Cells.Item($newRow, 10).Value = "Fly (Austin)"
$ws.Cells.Item($newRow, 11).Value = 2012
$ws.Cells.Item($newRow, 12).Value = "26/06/2012"
$ws.Cells.Item($newRow, 13).Value = "PMC3679285"
$ws.Cells.Item($newRow, 15).Value = "10.4161/fly.19695"

# Grow the table (ListObject) to include the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:O120"))
